$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1273.5
$ws.Range("I98").Value = 1205.3846
$ws.Range("J98").Value = 1568.6666
$ws.Range("K98").Value = 1205.3846
$ws.Range("L98").Value = 1568.6666
$ws.Range("M98").Value = 292.6153999999999
$ws.Range("N98").Value = -4564.6666

$ws.Range("H122").Value = 1273.5
$ws.Range("I122").Value = 1205.3846
$ws.Range("J122").Value = 1568.6666
$ws.Range("K122").Value = 3616.1538
$ws.Range("L122").Value = 4705.9998
$ws.Range("M122").Value = -1166.1538
$ws.Range("N122").Value = -9605.9998

$ws.Range("H132").Value = 1865
$ws.Range("I132").Value = 1863.8334
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 5591.5002
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -3061.5002
$ws.Range("N132").Value = -10760

$ws.Range("H135").Value = 55556532
$ws.Range("I135").Value = 21740172
$ws.Range("K135").Value = 195661548
$ws.Range("M135").Value = -195659013

$ws.Range("H137").Value = 2485.6072
$ws.Range("I137").Value = 2064.75
$ws.Range("K137").Value = 6194.25
$ws.Range("M137").Value = -3644.25

$ws.Range("H138").Value = 2662.8271
$ws.Range("I138").Value = 1038.2094
$ws.Range("J138").Value = 4501.2104
$ws.Range("K138").Value = 3114.6282
$ws.Range("L138").Value = 13503.6312
$ws.Range("M138").Value = 2025.3718
$ws.Range("N138").Value = -23783.6312

$ws.Range("H141").Value = 2494.647
$ws.Range("I141").Value = 1893.7241
$ws.Range("K141").Value = 5681.1723
$ws.Range("M141").Value = -501.1723000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18091.594
$ws.Range("I32").Value = 19099.598
$ws.Range("K32").Value = 19099.598
$ws.Range("M32").Value = -18812.598

$ws.Range("H45").Value = 1234.8387
$ws.Range("I45").Value = 1172.6863
$ws.Range("J45").Value = 1523
$ws.Range("K45").Value = 1172.6863
$ws.Range("L45").Value = 1523
$ws.Range("M45").Value = -795.6863000000001
$ws.Range("N45").Value = -2277

$ws.Range("H61").Value = 6414.051
$ws.Range("I61").Value = 3482.525
$ws.Range("J61").Value = 12585.685
$ws.Range("K61").Value = 3482.525
$ws.Range("L61").Value = 12585.685
$ws.Range("M61").Value = -3270.525
$ws.Range("N61").Value = -13009.685

$ws.Range("H74").Value = 6110.36
$ws.Range("I74").Value = 2229.3809
$ws.Range("J74").Value = 26485.5
$ws.Range("K74").Value = 2229.3809
$ws.Range("L74").Value = 26485.5
$ws.Range("M74").Value = -1355.3809
$ws.Range("N74").Value = -28233.5

$ws.Range("H77").Value = 6110.36
$ws.Range("I77").Value = 2229.3809
$ws.Range("J77").Value = 26485.5
$ws.Range("K77").Value = 11146.9045
$ws.Range("L77").Value = 132427.5
$ws.Range("M77").Value = -6778.904500000001
$ws.Range("N77").Value = -141163.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 1713.3334
$ws.Range("I132").Value = 1326.4445
$ws.Range("J132").Value = 2874
$ws.Range("K132").Value = 3979.3335
$ws.Range("L132").Value = 8622
$ws.Range("M132").Value = -1449.3335
$ws.Range("N132").Value = -13682

$ws.Range("H136").Value = 6414.051
$ws.Range("I136").Value = 3482.525
$ws.Range("J136").Value = 12585.685
$ws.Range("K136").Value = 10447.575
$ws.Range("L136").Value = 37757.055
$ws.Range("M136").Value = -7897.575000000001
$ws.Range("N136").Value = -42857.055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1994.6724
$ws.Range("I31").Value = 1398.7906
$ws.Range("J31").Value = 3702.8667
$ws.Range("K31").Value = 1398.7906
$ws.Range("L31").Value = 3702.8667
$ws.Range("M31").Value = -1103.7906
$ws.Range("N31").Value = -4292.8667

$ws.Range("H34").Value = 1994.6724
$ws.Range("I34").Value = 1398.7906
$ws.Range("J34").Value = 3702.8667
$ws.Range("K34").Value = 1398.7906
$ws.Range("L34").Value = 3702.8667
$ws.Range("M34").Value = -1196.7906
$ws.Range("N34").Value = -4106.8667

$ws.Range("H94").Value = 1605.6
$ws.Range("J94").Value = 1605.6
$ws.Range("L94").Value = 1605.6
$ws.Range("N94").Value = -2507.6

$ws.Range("H99").Value = 3539.2856
$ws.Range("I99").Value = 3102.2
$ws.Range("K99").Value = 3102.2
$ws.Range("M99").Value = -1604.2

$ws.Range("H126").Value = 3539.2856
$ws.Range("I126").Value = 3102.2
$ws.Range("K126").Value = 9306.599999999999
$ws.Range("M126").Value = -6836.599999999999

$ws.Range("H132").Value = 3116.4856
$ws.Range("I132").Value = 3431.375
$ws.Range("J132").Value = 2429.4546
$ws.Range("K132").Value = 10294.125
$ws.Range("L132").Value = 7288.3638
$ws.Range("M132").Value = -7764.125
$ws.Range("N132").Value = -12348.3638

$ws.Range("H134").Value = 2347.6128
$ws.Range("J134").Value = 3272.5715
$ws.Range("L134").Value = 9817.7145
$ws.Range("N134").Value = -14887.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 31250582
$ws.Range("I5").Value = 626.125
$ws.Range("J5").Value = 62500540
$ws.Range("K5").Value = 1878.375
$ws.Range("L5").Value = 187501620
$ws.Range("M5").Value = -1766.375
$ws.Range("N5").Value = -187501844

$ws.Range("H68").Value = 923.7857
$ws.Range("J68").Value = 966.63635
$ws.Range("L68").Value = 2899.90905
$ws.Range("N68").Value = -4521.90905

$ws.Range("H71").Value = 923.7857
$ws.Range("J71").Value = 966.63635
$ws.Range("L71").Value = 8699.727150000001
$ws.Range("N71").Value = -16811.72715

$ws.Range("H86").Value = 1787.4286
$ws.Range("I86").Value = 2141.8
$ws.Range("J86").Value = 901.5
$ws.Range("K86").Value = 6425.400000000001
$ws.Range("L86").Value = 2704.5
$ws.Range("M86").Value = -5239.400000000001
$ws.Range("N86").Value = -5076.5

$ws.Range("H89").Value = 1787.4286
$ws.Range("I89").Value = 2141.8
$ws.Range("J89").Value = 901.5
$ws.Range("K89").Value = 19276.2
$ws.Range("L89").Value = 8113.5
$ws.Range("M89").Value = -13348.2
$ws.Range("N89").Value = -19969.5

$ws.Range("H131").Value = 26180.71
$ws.Range("I131").Value = 1316
$ws.Range("J131").Value = 29948.092
$ws.Range("K131").Value = 3948
$ws.Range("L131").Value = 89844.276
$ws.Range("M131").Value = 1092
$ws.Range("N131").Value = -99924.276

$ws.Range("H135").Value = 31250582
$ws.Range("I135").Value = 626.125
$ws.Range("J135").Value = 62500540
$ws.Range("K135").Value = 5635.125
$ws.Range("L135").Value = 562504860
$ws.Range("M135").Value = -3100.125
$ws.Range("N135").Value = -562509930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 3191.7903
$ws.Range("I132").Value = 1581.02
$ws.Range("J132").Value = 9903.333000000001
$ws.Range("K132").Value = 4743.059999999999
$ws.Range("L132").Value = 29709.999
$ws.Range("M132").Value = -2213.059999999999
$ws.Range("N132").Value = -34769.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10050000
$ws.Range("I2").Value = 20000000
$ws.Range("K2").Value = 20000000
$ws.Range("M2").Value = -19999888

$ws.Range("H16").Value = 1175.4706
$ws.Range("I16").Value = 883.4167
$ws.Range("J16").Value = 1876.4
$ws.Range("K16").Value = 883.4167
$ws.Range("L16").Value = 1876.4
$ws.Range("M16").Value = -713.4167
$ws.Range("N16").Value = -2216.4

$ws.Range("H40").Value = 3518.7727
$ws.Range("I40").Value = 3493.4285
$ws.Range("J40").Value = 3563.125
$ws.Range("K40").Value = 3493.4285
$ws.Range("L40").Value = 3563.125
$ws.Range("M40").Value = -3357.4285
$ws.Range("N40").Value = -3835.125

$ws.Range("H46").Value = 1968.1
$ws.Range("I46").Value = 2000.25
$ws.Range("J46").Value = 1946.6666
$ws.Range("K46").Value = 2000.25
$ws.Range("L46").Value = 1946.6666
$ws.Range("M46").Value = -1812.25
$ws.Range("N46").Value = -2322.6666

$ws.Range("H55").Value = 611.7143
$ws.Range("I55").Value = 660
$ws.Range("K55").Value = 660
$ws.Range("M55").Value = -487

$ws.Range("H132").Value = 4015.6792
$ws.Range("I132").Value = 3982.257
$ws.Range("K132").Value = 11946.771
$ws.Range("M132").Value = -9416.771000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2449
$ws.Range("I81").Value = 1582.7
$ws.Range("J81").Value = 3115.3845
$ws.Range("K81").Value = 3165.4
$ws.Range("L81").Value = 6230.769
$ws.Range("M81").Value = -2104.4
$ws.Range("N81").Value = -8352.769

$ws.Range("H84").Value = 2449
$ws.Range("I84").Value = 1582.7
$ws.Range("J84").Value = 3115.3845
$ws.Range("K84").Value = 15827
$ws.Range("L84").Value = 31153.845
$ws.Range("M84").Value = -10523
$ws.Range("N84").Value = -41761.845

$ws.Range("H118").Value = 68900
$ws.Range("J118").Value = 68900
$ws.Range("L118").Value = 68900
$ws.Range("N118").Value = -72214

$ws.Range("H126").Value = 1710.2727
$ws.Range("I126").Value = 2703.7144
$ws.Range("J126").Value = 1246.6666
$ws.Range("K126").Value = 8111.1432
$ws.Range("L126").Value = 3739.9998
$ws.Range("M126").Value = -5641.1432
$ws.Range("N126").Value = -8679.9998

$ws.Range("H136").Value = 5650.1357
$ws.Range("I136").Value = 4706.129
$ws.Range("J136").Value = 6695.2856
$ws.Range("K136").Value = 14118.387
$ws.Range("L136").Value = 20085.8568
$ws.Range("M136").Value = -11568.387
$ws.Range("N136").Value = -25185.8568
